$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.85"
$ws.Range("E2").Value = "'-1.32%"
$ws.Range("D3").Value = "'40.85"
$ws.Range("E3").Value = "'-0.63%"
$ws.Range("D4").Value = "'5.042"
$ws.Range("E4").Value = "'-1.18%"
$ws.Range("D5").Value = "'0.07638"
$ws.Range("E5").Value = "'-3.06%"
$ws.Range("D6").Value = "'4.264"
$ws.Range("E6").Value = "'-1.54%"
$ws.Range("D7").Value = "'1.617"
$ws.Range("E7").Value = "'-4.26%"
$ws.Range("D9").Value = "'0.9094"
$ws.Range("E9").Value = "'-1.77%"
$ws.Range("D10").Value = "'0.1010"
$ws.Range("E10").Value = "'-6.72%"
$ws.Range("D11").Value = "'0.1770"
$ws.Range("E11").Value = "'-0.86%"
$ws.Range("D12").Value = "'0.09181"
$ws.Range("E12").Value = "'0.84%"
$ws.Range("D13").Value = "'0.04303"
$ws.Range("E13").Value = "'-2.27%"
$ws.Range("D14").Value = "'0.1054"
$ws.Range("E14").Value = "'-0.62%"
$ws.Range("D15").Value = "'0.001257"
$ws.Range("E15").Value = "'-0.79%"
$ws.Range("D16").Value = "'0.005834"
$ws.Range("E16").Value = "'-3.19%"
$ws.Range("E17").Value = "'-0.65%"
$ws.Range("E18").Value = "'-1.42%"
$ws.Range("D19").Value = "'6.781"
$ws.Range("E19").Value = "'-5.71%"
$ws.Range("D20").Value = "'0.1355"
$ws.Range("E20").Value = "'-1.49%"
$ws.Range("E21").Value = "'-2.80%"
$ws.Range("D22").Value = "'0.04161"
$ws.Range("E22").Value = "'0.10%"
$ws.Range("D23").Value = "'0.001225"
$ws.Range("E23").Value = "'-1.35%"
$ws.Range("D24").Value = "'0.004093"
$ws.Range("E24").Value = "'-1.43%"
$ws.Range("D25").Value = "'0.0001301"
$ws.Range("E25").Value = "'6.08%"
$ws.Range("D26").Value = "'0.0003009"
$ws.Range("E26").Value = "'0.51%"
$ws.Range("D38").Value = "'0.02412"
$ws.Range("E38").Value = "'-1.80%"
$ws.Range("D39").Value = "'0.05177"
$ws.Range("E39").Value = "'-2.83%"
$ws.Range("D40").Value = "'0.007767"
$ws.Range("E40").Value = "'-3.16%"
$ws.Range("D41").Value = "'0.1310"
$ws.Range("E41").Value = "'-3.47%"
$ws.Range("D42").Value = "'0.007089"
$ws.Range("E42").Value = "'-6.44%"
$ws.Range("D43").Value = "'0.001949"
$ws.Range("E43").Value = "'-2.40%"
$ws.Range("D44").Value = "'0.007462"
$ws.Range("E44").Value = "'-9.05%"
$ws.Range("D45").Value = "'0.3059"
$ws.Range("E45").Value = "'-1.45%"
$ws.Range("D46").Value = "'0.00006372"
$ws.Range("E46").Value = "'-5.77%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.42%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.004403"
$ws.Range("E48").Value = "'6.85%"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.005653"
$ws.Range("E49").Value = "'64.72%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.42%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.42%"
